# Auto-generated Excel COM-interop script to apply crypto data refresh
# Updates Price (D) and Volume(1h) (E) columns, and swaps three coin rows
# (Bittensor/dogwifhat, USDe/Stacks, Aave/OKB) that changed rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.278.64'
$ws.Range("E2").Value = '  +3.26%  '
$ws.Range("D3").Value = '3.207.66'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = '''218.45'
$ws.Range("E5").Value = '  +6.61%  '
$ws.Range("D6").Value = '''650.48'
$ws.Range("E6").Value = '  +7.12%  '
$ws.Range("D7").Value = '''0.397'
$ws.Range("E7").Value = '  +5.20%  '
$ws.Range("D8").Value = '''0.695'
$ws.Range("E8").Value = '  +5.46%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").Value = '3.209.31'
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("D11").Value = '''0.577'
$ws.Range("E11").Value = '  +8.41%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '''0.0000257'
$ws.Range("E13").Value = '  +5.82%  '
$ws.Range("D14").Value = '''5.42'
$ws.Range("E14").Value = '  +3.54%  '
$ws.Range("D15").Value = '''33.54'
$ws.Range("E15").Value = '  +5.14%  '
$ws.Range("D16").Value = '3.787.88'
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").Value = '89.716.20'
$ws.Range("E17").Value = '  +3.01%  '
$ws.Range("D18").Value = '3.199.36'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = '''3.45'
$ws.Range("E19").Value = '  +15.32%  '
$ws.Range("D20").Value = '''0.0000228'
$ws.Range("E20").Value = '  +76.75%  '
$ws.Range("D21").Value = '''13.58'
$ws.Range("E21").Value = '  +1.72%  '
$ws.Range("D22").Value = '''440.95'
$ws.Range("E22").Value = '  +6.74%  '
$ws.Range("D23").Value = '''8.70'
$ws.Range("E23").Value = '  +2.94%  '
$ws.Range("D24").Value = '''5.11'
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").Value = '''5.35'
$ws.Range("E25").Value = '  +3.85%  '
$ws.Range("D26").Value = '''12.06'
$ws.Range("E26").Value = '  +2.45%  '
$ws.Range("D27").Value = '''81.82'
$ws.Range("E27").Value = '  +11.69%  '
$ws.Range("D28").Value = '3.366.31'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '''0.160'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").Value = '''0.996'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").Value = '''4.09'
$ws.Range("E32").Value = '  +36.73%  '
$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").Value = '''550.96'
$ws.Range("E33").Value = '  +2.78%  '
$ws.Range("D34").Value = '''8.54'
$ws.Range("E34").Value = '  +3.54%  '
$ws.Range("D35").Value = '''7.12'
$ws.Range("E35").Value = '  +6.63%  '
$ws.Range("E36").Value = '  +6.08%  '
$ws.Range("E37").Value = '  +1.66%  '
$ws.Range("D38").Value = '''22.61'
$ws.Range("E38").Value = '  +3.62%  '
$ws.Range("D39").Value = '''22.41'
$ws.Range("E39").Value = '  +2.80%  '
$ws.Range("D40").Value = '''0.129'
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''1.95'
$ws.Range("E42").Value = '  +3.37%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = '''1.00'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").Value = '''0.377'
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("D45").Value = '''146.66'
$ws.Range("E45").Value = '  -0.49%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '''45.05'
$ws.Range("E46").Value = '  +4.33%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''174.99'
$ws.Range("E47").Value = '  +2.00%  '
$ws.Range("D48").Value = '''0.767'
$ws.Range("E48").Value = '  +10.75%  '
$ws.Range("D49").Value = '''0.125'
$ws.Range("E49").Value = '  -1.68%  '
$ws.Range("E50").Value = '  +0.87%  '
$ws.Range("D51").Value = '''0.626'
$ws.Range("E51").Value = '  +7.05%  '
